$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 628.8
$ws.Range("I34").Value = 628.8
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 628.8
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -425.8
$ws.Range("H36").Value = 628.8
$ws.Range("I36").Value = 628.8
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 628.8
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = 86.20000000000005
$ws.Range("H53").Value = 8902.352999999999
$ws.Range("I53").Value = 13688.363
$ws.Range("J53").Value = 128
$ws.Range("K53").Value = 13688.363
$ws.Range("L53").Value = 128
$ws.Range("M53").Value = -13051.363
$ws.Range("N53").Value = -1402
$ws.Range("H76").Value = 28574038
$ws.Range("I76").Value = 35716884
$ws.Range("J76").Value = 2661.4285
$ws.Range("K76").Value = 35716884
$ws.Range("L76").Value = 2661.4285
$ws.Range("M76").Value = -35716569
$ws.Range("N76").Value = -3291.4285
$ws.Range("H79").Value = 28574038
$ws.Range("I79").Value = 35716884
$ws.Range("J79").Value = 2661.4285
$ws.Range("K79").Value = 35716884
$ws.Range("L79").Value = 2661.4285
$ws.Range("M79").Value = -35715792
$ws.Range("N79").Value = -4845.4285
$ws.Range("H116").Value = 19850354
$ws.Range("I116").Value = 11907930
$ws.Range("J116").Value = 27792778
$ws.Range("K116").Value = 11907930
$ws.Range("L116").Value = 27792778
$ws.Range("M116").Value = -11904488
$ws.Range("N116").Value = -27799662
$ws.Range("H128").Value = 15743.75
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 15743.75
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 15743.75
$ws.Range("N128").Value = -25703.75
$ws.Range("H132").Value = 3028452.5
$ws.Range("I132").Value = 628691.2
$ws.Range("J132").Value = 22226542
$ws.Range("K132").Value = 1886073.6
$ws.Range("L132").Value = 66679626
$ws.Range("M132").Value = -1883543.6
$ws.Range("N132").Value = -66684686

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2440.923
$ws.Range("I110").Value = 951.8889
$ws.Range("J110").Value = 5791.25
$ws.Range("K110").Value = 951.8889
$ws.Range("L110").Value = 5791.25
$ws.Range("M110").Value = 1093.1111
$ws.Range("N110").Value = -9881.25
$ws.Range("H132").Value = 33732304
$ws.Range("I132").Value = 39302916
$ws.Range("J132").Value = 15627815
$ws.Range("K132").Value = 117908748
$ws.Range("L132").Value = 46883445
$ws.Range("M132").Value = -117906218
$ws.Range("N132").Value = -46888505

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 667681.8
$ws.Range("I107").Value = 770293.9399999999
$ws.Range("J107").Value = 703
$ws.Range("K107").Value = 770293.9399999999
$ws.Range("L107").Value = 703
$ws.Range("M107").Value = -768373.9399999999
$ws.Range("N107").Value = -4543
$ws.Range("H134").Value = 12080029
$ws.Range("I134").Value = 16216071
$ws.Range("J134").Value = 47905.727
$ws.Range("K134").Value = 48648213
$ws.Range("L134").Value = 143717.181
$ws.Range("M134").Value = -48645678
$ws.Range("N134").Value = -148787.181

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 59999
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 59999
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 59999
$ws.Range("N20").Value = -60471
$ws.Range("H30").Value = 59999
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 59999
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 59999
$ws.Range("N30").Value = -60181
$ws.Range("H31").Value = 3090824.8
$ws.Range("I31").Value = 7577240.5
$ws.Range("J31").Value = 6413.625
$ws.Range("K31").Value = 7577240.5
$ws.Range("L31").Value = 6413.625
$ws.Range("M31").Value = -7576945.5
$ws.Range("N31").Value = -7003.625
$ws.Range("H34").Value = 3090824.8
$ws.Range("I34").Value = 7577240.5
$ws.Range("J34").Value = 6413.625
$ws.Range("K34").Value = 7577240.5
$ws.Range("L34").Value = 6413.625
$ws.Range("M34").Value = -7577038.5
$ws.Range("N34").Value = -6817.625
$ws.Range("H58").Value = 2681176.2
$ws.Range("I58").Value = 13657.625
$ws.Range("J58").Value = 5052304
$ws.Range("K58").Value = 13657.625
$ws.Range("L58").Value = 5052304
$ws.Range("M58").Value = -13454.625
$ws.Range("N58").Value = -5052710
$ws.Range("H128").Value = 59999
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 59999
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 59999
$ws.Range("N128").Value = -69959
$ws.Range("H132").Value = 1850.5483
$ws.Range("I132").Value = 1323.1305
$ws.Range("J132").Value = 3366.875
$ws.Range("K132").Value = 3969.3915
$ws.Range("L132").Value = 10100.625
$ws.Range("M132").Value = -1439.3915
$ws.Range("N132").Value = -15160.625
$ws.Range("H134").Value = 1484624.1
$ws.Range("I134").Value = 2418.125
$ws.Range("J134").Value = 3640560.2
$ws.Range("K134").Value = 7254.375
$ws.Range("L134").Value = 10921680.6
$ws.Range("M134").Value = -4719.375
$ws.Range("N134").Value = -10926750.6
$ws.Range("H136").Value = 2681176.2
$ws.Range("I136").Value = 13657.625
$ws.Range("J136").Value = 5052304
$ws.Range("K136").Value = 40972.875
$ws.Range("L136").Value = 15156912
$ws.Range("M136").Value = -38422.875
$ws.Range("N136").Value = -15162012

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 855450.75
$ws.Range("I107").Value = 2136889.8
$ws.Range("J107").Value = 1158
$ws.Range("K107").Value = 6410669.399999999
$ws.Range("L107").Value = 3474
$ws.Range("M107").Value = -6408749.399999999
$ws.Range("N107").Value = -7314
$ws.Range("H109").Value = 5245.3076
$ws.Range("I109").Value = 1599.75
$ws.Range("J109").Value = 6865.5557
$ws.Range("K109").Value = 4799.25
$ws.Range("L109").Value = 20596.6671
$ws.Range("M109").Value = -3759.25
$ws.Range("N109").Value = -22676.6671
$ws.Range("H129").Value = 1953.3235
$ws.Range("I129").Value = 699.8421
$ws.Range("J129").Value = 3541.0667
$ws.Range("K129").Value = 2099.5263
$ws.Range("L129").Value = 10623.2001
$ws.Range("M129").Value = 2900.4737
$ws.Range("N129").Value = -20623.2001
$ws.Range("H131").Value = 893.17975
$ws.Range("I131").Value = 311.1111
$ws.Range("J131").Value = 958.6625
$ws.Range("K131").Value = 933.3333
$ws.Range("L131").Value = 2875.9875
$ws.Range("M131").Value = 4106.6667
$ws.Range("N131").Value = -12955.9875
$ws.Range("H134").Value = 2866.6667
$ws.Range("I134").Value = 2300
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 6900
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -1830
$ws.Range("N134").Value = -22140
$ws.Range("H139").Value = 24024.682
$ws.Range("I139").Value = 28919.166
$ws.Range("J139").Value = 1999.5
$ws.Range("K139").Value = 86757.49800000001
$ws.Range("L139").Value = 5998.5
$ws.Range("M139").Value = -81617.49800000001
$ws.Range("N139").Value = -16278.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 35165.75
$ws.Range("I113").Value = 1250
$ws.Range("J113").Value = 46471
$ws.Range("K113").Value = 1250
$ws.Range("L113").Value = 46471
$ws.Range("M113").Value = 920
$ws.Range("N113").Value = -50811
$ws.Range("H126").Value = 5892.2856
$ws.Range("I126").Value = 9128.923000000001
$ws.Range("J126").Value = 3087.2
$ws.Range("K126").Value = 27386.769
$ws.Range("L126").Value = 9261.599999999999
$ws.Range("M126").Value = -24916.769
$ws.Range("N126").Value = -14201.6
$ws.Range("H132").Value = 6910460
$ws.Range("I132").Value = 7989046
$ws.Range("J132").Value = 5052894.5
$ws.Range("K132").Value = 23967138
$ws.Range("L132").Value = 15158683.5
$ws.Range("M132").Value = -23964608
$ws.Range("N132").Value = -15163743.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 35715096
$ws.Range("I22").Value = 1125
$ws.Range("J22").Value = 50000684
$ws.Range("K22").Value = 1125
$ws.Range("L22").Value = 50000684
$ws.Range("M22").Value = -830
$ws.Range("N22").Value = -50001274
$ws.Range("H27").Value = 35715096
$ws.Range("I27").Value = 1125
$ws.Range("J27").Value = 50000684
$ws.Range("K27").Value = 1125
$ws.Range("L27").Value = 50000684
$ws.Range("M27").Value = -1018
$ws.Range("N27").Value = -50000898

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 7879.4
$ws.Range("I39").Value = 7777
$ws.Range("J39").Value = 7905
$ws.Range("K39").Value = 7777
$ws.Range("L39").Value = 7905
$ws.Range("M39").Value = -7364
$ws.Range("N39").Value = -8731
$ws.Range("H42").Value = 8830
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 8830
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 8830
$ws.Range("N42").Value = -9586
$ws.Range("H43").Value = 7490
$ws.Range("I43").Value = 5000
$ws.Range("J43").Value = 7905
$ws.Range("K43").Value = 5000
$ws.Range("L43").Value = 7905
$ws.Range("M43").Value = -4851
$ws.Range("N43").Value = -8203
$ws.Range("H132").Value = 692273.2
$ws.Range("I132").Value = 1864.5161
$ws.Range("J132").Value = 6807321.5
$ws.Range("K132").Value = 5593.5483
$ws.Range("L132").Value = 20421964.5
$ws.Range("M132").Value = -3063.5483
$ws.Range("N132").Value = -20427024.5
